# Remove the "configuracao_hard" (F) and "configuracao_soft" (G) columns.
# Deleting entire columns shifts everything to the right of them
# (funcao_api, qt_bytes, qt_requisicoes, time_stamp_init, time_stamp_fin,
# "latencia ( J - I )") two positions to the left, along with their
# column widths/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1:G1").EntireColumn.Delete()

# Update the selection to match the target state.
$ws.Range("M11").Select()
